$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("[class]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$classPara = $rng.Paragraphs(1)
$classParaRange = $classPara.Range
$endBefore = $classParaRange.End
Write-Host "endBefore=$endBefore"

$insertPoint = $d.Range($endBefore, $endBefore)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"></pkg:package>'

$ooxml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">ng g interceptor error </w:t></w:r><w:r><w:t>--</w:t></w:r><w:r><w:t>skip-tests</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>[interceptor]</w:t></w:r></w:p>
'@

try {
    $insertPoint.InsertXML($ooxml)
    Write-Host "InsertXML ok"
} catch {
    Write-Host "InsertXML failed: $_"
}
